$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (ColumnWidth units are offset ~0.9 from the stored XML "width") ---
$ws.Range("A:A").ColumnWidth = 24.1
$ws.Range("B:B").ColumnWidth = 24.1
$ws.Range("C:C").ColumnWidth = 22.1
$ws.Range("D:D").ColumnWidth = 29.1
$ws.Range("E:E").ColumnWidth = 29.1

# --- Text ("@") number format for the Phone Number (C) and Password (E) columns ---
$ws.Range("C1:C4").NumberFormat = "@"
$ws.Range("E1:E4").NumberFormat = "@"

# --- Cell values, written in an order that reproduces the target shared-string table ---
$ws.Range("A2").Value = "mohamed.h.eladwy"
$ws.Range("B2").Value = "Mohamed Eladwy"
$ws.Range("B1").Value = "Full Name:"
$ws.Range("C1").Value = "Phone Number:"
$ws.Range("D1").Value = "Email:"
$ws.Range("E1").Value = "Password:"
$ws.Range("A1").Value = "ID:"
$ws.Range("D2").Value = "mohamed.h.eladwy@gmail.com"
$ws.Range("A3").Value = "Yusuf.Bdr132"
$ws.Range("B3").Value = "Yusuf Elsayad Bdr"
$ws.Range("C3").Value = "01142340941"
$ws.Range("D3").Value = "yusuf.bdr.123@gmail.com"
$ws.Range("A4").Value = "Amr.Elsayed.Elhenawy"
$ws.Range("B4").Value = "Amr Elsayed Elhenawy"
$ws.Range("C4").Value = "01121753452"
$ws.Range("D4").Value = "amr.elhenawy@gmail.com"
$ws.Range("E4").Value = "56/89/2003#Amr"
$ws.Range("E3").Value = "15/7/2006#Yusuf"
$ws.Range("E2").Value = "132@Hussein"
$ws.Range("C2").Value = "01120664373"

# --- Hyperlinks on the Email column for the 3 data rows ---
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:mohamed.h.eladwy@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:yusuf.bdr.123@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:amr.elhenawy@gmail.com")

# --- E2 carries the Hyperlink look (underline / theme color) plus the Text number format ---
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "132@Hussein"

# --- Page setup / selection ---
$ws.PageSetup.Orientation = 1
[void]$ws.Range("E8").Select()

Write-Host "done"
